$p = $ppt.ActivePresentation
Write-Output "HasNotesMaster: $($p.HasNotesMaster)"
